# Monthly rollover update for LOZANO MOLINA TITO workbook.
#
# "VENTA MENSUAL" sheet: the tracked months shift forward by one
# (julio/agosto/septiembre/octubre -> agosto/septiembre/octubre/noviembre),
# so every data row's values shift one column to the left (new C = old D,
# new D = old E, new E = old F) and the newest month column (F, "noviembre")
# starts out empty (0). Column widths are re-fitted to the new header text.
#
# "VENTAS POR GRUPO" sheet: the figures that belonged to the month that
# rolled out of the "VENTA MENSUAL" window are cleared back to 0, and the
# "X de 30" completion counters in row 32 are decremented to match.

$wb = $excel.ActiveWorkbook

# ---- VENTA MENSUAL ---------------------------------------------------
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Shift header labels one month forward.
$ws.Range("C1").Value = "agosto"
$ws.Range("D1").Value = "septiembre"
$ws.Range("E1").Value = "octubre"
$ws.Range("F1").Value = "noviembre"

# Shift each data row's monthly figures one column to the left; the new
# (rightmost, newest) month has no data yet.
for ($r = 2; $r -le 32; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2

    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $fVal
    $ws.Cells.Item($r, 6).Value = 0
}

# Re-fit the column widths for the new header text (observed offset of
# 5/6 character between the ColumnWidth API and the stored OOXML width).
$ws.Columns.Item(4).ColumnWidth = 16 - 0.8333333333
$ws.Columns.Item(5).ColumnWidth = 13 - 0.8333333333
$ws.Columns.Item(6).ColumnWidth = 15 - 0.8333333333

# ---- VENTAS POR GRUPO --------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Figures for the month that rolled out of the tracking window are reset.
$ws2.Range("L4").Value = 0
$ws2.Range("D6").Value = 0
$ws2.Range("M16").Value = 0
$ws2.Range("D20").Value = 0
$ws2.Range("D22").Value = 0
$ws2.Range("D24").Value = 0

# Completion counters on the totals row (32) drop to 0 since all of their
# contributing entries were just cleared.
$ws2.Range("D32").Value = "0 de 30"
$ws2.Range("L32").Value = "0 de 30"
$ws2.Range("M32").Value = "0 de 30"
